$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in row 5 (P5, Q5, R5) ---
$ws.Range("P5").Value = 4.4000000000000004
$ws.Range("Q5").Value = 2.9
$ws.Range("R5").Value = 3.2

# --- Add new column S (year 2022) mirroring column R's formatting ---
# Row 4 header cell (style s="7" in the target XML)
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Row 5 data cell (style s="15" in the target XML)
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 3.4

# --- Move the active selection to T4 ---
$ws.Range("T4").Select()
